# Update the "time_taken" timestamps on the existing "data" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "2021-10-05 14:20:16.260096"
$ws.Range("F3").Value = "2021-10-05 14:20:16.260104"
$ws.Range("F4").Value = "2021-10-05 14:20:16.260107"
$ws.Range("F5").Value = "2021-10-05 14:20:16.260110"
$ws.Range("F6").Value = "2021-10-05 14:20:16.260113"
$ws.Range("F7").Value = "2021-10-05 14:20:16.260116"
$ws.Range("F8").Value = "2021-10-05 14:20:16.260119"
$ws.Range("F9").Value = "2021-10-05 14:20:16.260121"
$ws.Range("F10").Value = "2021-10-05 14:20:16.260124"
$ws.Range("F11").Value = "2021-10-05 14:20:16.260127"

# Add a new "metadata" worksheet right after the "data" sheet.
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Copy the header formatting (bold / centered / bordered) from the data
# sheet's header row onto the new metadata header row.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Familial hyperparathyroidism"
$meta.Range("C2").Value = 480
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.14"
$meta.Range("E2").Value = "2021-07-28T13:54:56.213879Z"
$meta.Range("F2").Value = "2021-10-05 14:20:16.256405"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/480/?format=json"

# Keep the "data" sheet as the active tab (matches original workbook view).
$ws.Activate()
$null = $ws.Range("A1").Select()
